# The sheet holds a year-indexed table (A: year label, B..J: case counts).
# The edit drops the three oldest years (2007-2009, rows 2-4) and appends a
# new year (2021, row 13) that only has a total in column J, mirroring the
# 2019/2020 rows that already existed with blank B..I cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 2007/2008/2009 rows; everything below shifts up by 3.
$ws.Rows("2:4").Delete()

# After the shift, the last data row is 2020 (row 12). Add 2021 as row 13,
# copying the 2020 row's formatting (bold/centered/bordered label style)
# onto the new label cell so it matches the rest of column A.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"

# Columns B..I have no data yet for 2021 (left blank, same as 2019/2020).
$ws.Range("B13:I13").Value = ""

# Only the overall total (column J) is known for 2021.
$ws.Range("J13").Value = 6380
